# Generate Report for Handoff
# Update status text + timestamps across the Overview / zh-cn / de-de sheets,
# and tighten the "Status"/date columns' width to match the new shorter text.

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn     = $wb.Worksheets.Item("zh-cn")
$dede     = $wb.Worksheets.Item("de-de")

# --- Text / timestamp updates -------------------------------------------------

# Status column: "Handed back: in sync with en-US" -> "Ready for handoff"
$overview.Range("E2").Value = "Ready for handoff"
$overview.Range("F2").Value = "Ready for handoff"
$zhcn.Range("C2").Value = "Ready for handoff"
$dede.Range("C2").Value = "Ready for handoff"

# Latest HO Xliff Generate Date / Latest Handback DateTime refresh
$overview.Range("G2").Value = "2016-08-30 13:04:16"
$dede.Range("H2").Value = "2016-08-30 13:04:16"

# Latest Handoff Datetime refresh (zh-cn)
$zhcn.Range("H2").Value = "2016-08-30 13:04:02"

# --- Column width updates -----------------------------------------------------
# Target stored width is 17.2159881591797 "characters"; Excel's ColumnWidth
# setter always snaps to a whole-pixel grid (standard width<->pixel rounding),
# so feed it the character width whose pixel-grid result lands nearest the
# target (16.3333... -> stored 17.1666...).

$overview.Columns.Item(5).ColumnWidth = 16.3333333333333
$overview.Columns.Item(6).ColumnWidth = 16.3333333333333
$zhcn.Columns.Item(3).ColumnWidth = 16.3333333333333
$dede.Columns.Item(3).ColumnWidth = 16.3333333333333
